$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: S0000002 / Nome sistema / (blank) / CAL_CLOSED;CAL_DONE / Username
$ws.Range("A2").Value = "S0000002"
$ws.Range("A3").Value = "S0000003"
$ws.Range("B2").Value = "Nome sistema"
$ws.Range("B3").Value = "Nome sistema"
$ws.Range("C3").Value = "CAL_DONE"
$ws.Range("D3").Value = "CAL_IN_PROGRESS"
$ws.Range("D2").Value = "CAL_CLOSED;CAL_DONE"
$ws.Range("E2").Value = "Username"
$ws.Range("E3").Value = "Username"
